# Met à jour les résultats
# Updates the evaluation-results table (second block, rows 17-24) with new
# measured values, adds an explanatory comment string, and moves the
# selection to the updated "Moyennes" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 18 (DayTripper / CD) -------------------------------------------
$ws.Range("H18").Value = 0.5555
$ws.Range("I18").Value = 0.4444
$ws.Range("J18").Value = 0.3571
$ws.Range("K18").Value = 0.4444
$ws.Range("L18").Value = 0.9444

# --- Row 20 (Aller-Retour diatonique / Synthétique) ----------------------
$ws.Range("J20").Value = 0.966

# --- Row 22 (Seven Nation Army / Enregistré) ------------------------------
$ws.Range("H22").Value = 0.825
$ws.Range("I22").Value = 0.675
$ws.Range("L22").Value = 0.8913

# --- Row 23 (Hardest Button to Button / Enregistré) -----------------------
$ws.Range("H23").Value = 0.7887
$ws.Range("I23").Value = 0.5774
$ws.Range("K23").Value = 0.1126
$ws.Range("L23").Value = 0.7605

# Row 24 averages recalc automatically (AVERAGE formulas over H18:H23, etc.)

# --- New "Changements importants" note ------------------------------------
$ws.Range("H35").Value = "Conséquence => l'onset est en retard par rapport à la réalité => la note est mal analysée en ton"

# --- Update current selection ---------------------------------------------
$ws.Range("L24").Select()
